$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated "Price" (column D) and "Volume(1h)" (column E) cells as scraped
# by the cryptos-list GitHub Action run.
$updates = @(
    @{ Cell = "D2"; Value = '30.227.46'; ForceText = $false },
    @{ Cell = "E2"; Value = '  +0.45%  '; ForceText = $false },
    @{ Cell = "D3"; Value = '1.990.69'; ForceText = $false },
    @{ Cell = "E3"; Value = '  +5.99%  '; ForceText = $false },
    @{ Cell = "E4"; Value = '  -0.27%  '; ForceText = $false },
    @{ Cell = "D5"; Value = '324.52'; ForceText = $true },
    @{ Cell = "E5"; Value = '  +1.56%  '; ForceText = $false },
    @{ Cell = "E6"; Value = '  -0.26%  '; ForceText = $false },
    @{ Cell = "D7"; Value = '0.5096'; ForceText = $true },
    @{ Cell = "E7"; Value = '  +1.15%  '; ForceText = $false },
    @{ Cell = "D8"; Value = '0.4129'; ForceText = $true },
    @{ Cell = "E8"; Value = '  +4.34%  '; ForceText = $false },
    @{ Cell = "D9"; Value = '0.08665'; ForceText = $true },
    @{ Cell = "E9"; Value = '  +5.46%  '; ForceText = $false },
    @{ Cell = "D10"; Value = '1.130'; ForceText = $true },
    @{ Cell = "E10"; Value = '  +3.39%  '; ForceText = $false },
    @{ Cell = "D11"; Value = '42.72'; ForceText = $true },
    @{ Cell = "E11"; Value = '  +1.46%  '; ForceText = $false },
    @{ Cell = "D12"; Value = '24.38'; ForceText = $true },
    @{ Cell = "E12"; Value = '  +3.29%  '; ForceText = $false },
    @{ Cell = "D13"; Value = '1.992.51'; ForceText = $false },
    @{ Cell = "E13"; Value = '  +5.48%  '; ForceText = $false },
    @{ Cell = "D14"; Value = '6.478'; ForceText = $true },
    @{ Cell = "E14"; Value = '  +2.87%  '; ForceText = $false },
    @{ Cell = "D15"; Value = '7.381'; ForceText = $true },
    @{ Cell = "E15"; Value = '  +2.58%  '; ForceText = $false },
    @{ Cell = "E16"; Value = '  -0.29%  '; ForceText = $false },
    @{ Cell = "E17"; Value = '  +2.41%  '; ForceText = $false },
    @{ Cell = "E18"; Value = '  +2.27%  '; ForceText = $false },
    @{ Cell = "D19"; Value = '0.06548'; ForceText = $true },
    @{ Cell = "E19"; Value = '  +1.31%  '; ForceText = $false },
    @{ Cell = "D20"; Value = '18.74'; ForceText = $true },
    @{ Cell = "E20"; Value = '  +3.69%  '; ForceText = $false },
    @{ Cell = "E21"; Value = '  -0.25%  '; ForceText = $false },
    @{ Cell = "D22"; Value = '6.070'; ForceText = $true },
    @{ Cell = "E22"; Value = '  +4.03%  '; ForceText = $false },
    @{ Cell = "D23"; Value = '30.295.32'; ForceText = $false },
    @{ Cell = "E23"; Value = '  +0.68%  '; ForceText = $false },
    @{ Cell = "D24"; Value = '11.57'; ForceText = $true },
    @{ Cell = "E24"; Value = '  +3.61%  '; ForceText = $false },
    @{ Cell = "D25"; Value = '2.202'; ForceText = $true },
    @{ Cell = "E25"; Value = '  +1.42%  '; ForceText = $false },
    @{ Cell = "D26"; Value = '2.225.26'; ForceText = $false },
    @{ Cell = "E26"; Value = '  +6.09%  '; ForceText = $false },
    @{ Cell = "D27"; Value = '22.50'; ForceText = $true },
    @{ Cell = "E27"; Value = '  +6.20%  '; ForceText = $false },
    @{ Cell = "D28"; Value = '163.02'; ForceText = $true },
    @{ Cell = "E28"; Value = '  +1.26%  '; ForceText = $false },
    @{ Cell = "D29"; Value = '2.350'; ForceText = $true },
    @{ Cell = "E29"; Value = '  +4.56%  '; ForceText = $false },
    @{ Cell = "D30"; Value = '130.49'; ForceText = $true },
    @{ Cell = "E30"; Value = '  +2.48%  '; ForceText = $false },
    @{ Cell = "E31"; Value = '  +5.20%  '; ForceText = $false },
    @{ Cell = "E32"; Value = '  +1.60%  '; ForceText = $false },
    @{ Cell = "D33"; Value = '6.055'; ForceText = $true },
    @{ Cell = "E33"; Value = '  +2.11%  '; ForceText = $false },
    @{ Cell = "D34"; Value = '3.815'; ForceText = $true },
    @{ Cell = "E34"; Value = '  +3.35%  '; ForceText = $false },
    @{ Cell = "D35"; Value = '1.308'; ForceText = $true },
    @{ Cell = "E35"; Value = '  +11.51%  '; ForceText = $false },
    @{ Cell = "D36"; Value = '0.02477'; ForceText = $true },
    @{ Cell = "E36"; Value = '  +2.21%  '; ForceText = $false },
    @{ Cell = "D37"; Value = '5.377'; ForceText = $true },
    @{ Cell = "E37"; Value = '  +1.67%  '; ForceText = $false },
    @{ Cell = "D38"; Value = '0.06517'; ForceText = $true },
    @{ Cell = "E38"; Value = '  +2.65%  '; ForceText = $false },
    @{ Cell = "D39"; Value = '0.2188'; ForceText = $true },
    @{ Cell = "E39"; Value = '  +2.64%  '; ForceText = $false },
    @{ Cell = "D40"; Value = '8.912'; ForceText = $true },
    @{ Cell = "E40"; Value = '  +4.89%  '; ForceText = $false },
    @{ Cell = "D41"; Value = '0.6571'; ForceText = $true },
    @{ Cell = "D42"; Value = '11.82'; ForceText = $true },
    @{ Cell = "E42"; Value = '  +4.75%  '; ForceText = $false },
    @{ Cell = "D43"; Value = '1.224'; ForceText = $true },
    @{ Cell = "E43"; Value = '  +1.04%  '; ForceText = $false },
    @{ Cell = "D44"; Value = '13.61'; ForceText = $true },
    @{ Cell = "E44"; Value = '  +4.06%  '; ForceText = $false },
    @{ Cell = "D45"; Value = '0.6100'; ForceText = $true },
    @{ Cell = "E45"; Value = '  +3.31%  '; ForceText = $false },
    @{ Cell = "D46"; Value = '2.194'; ForceText = $true },
    @{ Cell = "E46"; Value = '  +4.87%  '; ForceText = $false },
    @{ Cell = "D47"; Value = '3.660'; ForceText = $true },
    @{ Cell = "E47"; Value = '  +0.93%  '; ForceText = $false },
    @{ Cell = "D48"; Value = '124.26'; ForceText = $true },
    @{ Cell = "E48"; Value = '  +1.73%  '; ForceText = $false },
    @{ Cell = "E49"; Value = '  +1.18%  '; ForceText = $false },
    @{ Cell = "D50"; Value = '79.22'; ForceText = $true },
    @{ Cell = "E50"; Value = '  +2.32%  '; ForceText = $false },
    @{ Cell = "D51"; Value = '0.06863'; ForceText = $true },
    @{ Cell = "E51"; Value = '  +1.84%  '; ForceText = $false }
)

foreach ($u in $updates) {
    $cell = $ws.Range($u.Cell)
    if ($u.ForceText) {
        # Several prices (e.g. "1.130", "22.50") look like numbers with a
        # trailing zero; force the cell to Text format first so Excel
        # doesn't silently reinterpret/round them as numeric values.
        $cell.NumberFormat = "@"
        $cell.Value = $u.Value
        $cell.Style = "Normal"
    } else {
        $cell.Value = $u.Value
    }
}
